$wb = $excel.ActiveWorkbook
$cs = $wb.Sheets.Item("Burndown Chart")
Write-Host "Type:" $cs.Type
$members = $cs | Get-Member -Name "Type"
Write-Host $members
